$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 112.501404
$ws.Range("L2").Value = 7.947767
$ws.Range("R2").Value = 31.931187
$ws.Range("X2").Value = 102.775753
$ws.Range("Z2").Value = 72.080811
$ws.Range("AC2").Value = 50.728452

# Row 4 updates
$ws.Range("X4").Value = 118.145268
$ws.Range("AC4").Value = 47.621745

# Row 5 updates
$ws.Range("R5").Value = 43.050717

# Row 7 full rewrite (B7:AE7)
$ws.Range("B7").Value = 80.943821
$ws.Range("C7").Value = 134.507798
$ws.Range("D7").Value = 145.166512
$ws.Range("E7").Value = 121.602272
$ws.Range("F7").Value = 83.531475
$ws.Range("G7").Value = 30.315319
$ws.Range("H7").Value = 37.209151
$ws.Range("I7").Value = 36.880951
$ws.Range("J7").Value = 64.403607
$ws.Range("K7").Value = 55.711391
$ws.Range("L7").Value = 15.498488
$ws.Range("M7").Value = 32.673904
$ws.Range("N7").Value = 73.142105
$ws.Range("O7").Value = 68.236695
$ws.Range("P7").Value = 28.833284
$ws.Range("Q7").Value = 28.814277
$ws.Range("R7").Value = 44.043698
$ws.Range("S7").Value = 89.722458
$ws.Range("T7").Value = 64.403607
$ws.Range("U7").Value = 55.711391
$ws.Range("V7").Value = 41.598244
$ws.Range("W7").Value = 24.20197
$ws.Range("X7").Value = 145.872508
$ws.Range("Y7").Value = 132.509443
$ws.Range("Z7").Value = 71.147936
$ws.Range("AA7").Value = 22.222866
$ws.Range("AB7").Value = 152.690972
$ws.Range("AC7").Value = 116.703219
$ws.Range("AD7").Value = 64.403607
$ws.Range("AE7").Value = 55.711391
